# ============================================================================
# Add 2022-Q3 data:
#   1. "总计" (summary) sheet gets a new row-2 for the 2022-Q3 totals, pushing
#      every later quarter down by one row.
#   2. A brand-new "2022-Q3" worksheet (fund holdings detail) is inserted right
#      after "总计", pushing every later quarter sheet one tab to the right.
# ============================================================================

$wb = $excel.ActiveWorkbook
$total = $wb.Worksheets.Item(1)

# --- Step 1: "总计" sheet -- insert a new row 2 for the 2022-Q3 summary figures ---
$total.Rows.Item(2).Insert()
$total.Range("A3").Copy($total.Range("A2"))    # carry over the bordered/bold "index" style (s=2)
$total.Range("B2:D2").ClearFormats()            # Insert() inherited row-1 header formatting; strip it
$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 24
$total.Range("D2").Value = 2.37

# --- Step 2: insert the brand-new "2022-Q3" worksheet right after "总计" ---
$newSheet = $wb.Worksheets.Add($null, $total)
$newSheet.Name = "2022-Q3"

# Seed header-row + index-column styling by copying already-formatted cells from "总计"
$total.Range("B1").Copy($newSheet.Range("B1"))
$newSheet.Range("B1").Copy($newSheet.Range("C1"))
$newSheet.Range("B1").Copy($newSheet.Range("D1"))
$newSheet.Range("B1").Copy($newSheet.Range("E1"))
$newSheet.Range("B1").Copy($newSheet.Range("F1"))
$newSheet.Range("B1").Copy($newSheet.Range("G1"))
$newSheet.Range("B1").Copy($newSheet.Range("H1"))
$total.Range("A2").Copy($newSheet.Range("A2"))
$newSheet.Range("A2").Copy($newSheet.Range("A3:A25"))

# --- Step 3: header row text ---
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# --- Step 4: text-format data columns B:G so numeric-looking strings (fund codes,
#     sizes, ratios) keep their exact text, e.g. "008131", "26.10" -- column H (rank) stays numeric ---
$newSheet.Range("B2:G25").NumberFormat = "@"

# --- Step 5: data rows ---
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "008131"
$newSheet.Range("C2").Value = "景顺长城竞争优势混合"
$newSheet.Range("D2").Value = "26.10"
$newSheet.Range("E2").Value = "90.08"
$newSheet.Range("F2").Value = "2.68"
$newSheet.Range("G2").Value = "0.6995"
$newSheet.Range("H2").Value = 10
$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "519011"
$newSheet.Range("C3").Value = "海富通精选混合"
$newSheet.Range("D3").Value = "10.10"
$newSheet.Range("E3").Value = "77.01"
$newSheet.Range("F3").Value = "2.85"
$newSheet.Range("G3").Value = "0.2878"
$newSheet.Range("H3").Value = 10
$newSheet.Range("A4").Value = 2
$newSheet.Range("B4").Value = "260103"
$newSheet.Range("C4").Value = "景顺长城动力平衡混合"
$newSheet.Range("D4").Value = "10.57"
$newSheet.Range("E4").Value = "68.77"
$newSheet.Range("F4").Value = "2.43"
$newSheet.Range("G4").Value = "0.2569"
$newSheet.Range("H4").Value = 10
$newSheet.Range("A5").Value = 3
$newSheet.Range("B5").Value = "010220"
$newSheet.Range("C5").Value = "海富通消费核心资产混合A"
$newSheet.Range("D5").Value = "4.22"
$newSheet.Range("E5").Value = "91.44"
$newSheet.Range("F5").Value = "4.25"
$newSheet.Range("G5").Value = "0.1794"
$newSheet.Range("H5").Value = 9
$newSheet.Range("A6").Value = 4
$newSheet.Range("B6").Value = "001551"
$newSheet.Range("C6").Value = "天弘中证医药100指数型发起式 C"
$newSheet.Range("D6").Value = "8.58"
$newSheet.Range("E6").Value = "95.24"
$newSheet.Range("F6").Value = "1.51"
$newSheet.Range("G6").Value = "0.1296"
$newSheet.Range("H6").Value = 2
$newSheet.Range("A7").Value = 5
$newSheet.Range("B7").Value = "519015"
$newSheet.Range("C7").Value = "海富通精选贰号混合"
$newSheet.Range("D7").Value = "3.77"
$newSheet.Range("E7").Value = "77.19"
$newSheet.Range("F7").Value = "2.86"
$newSheet.Range("G7").Value = "0.1078"
$newSheet.Range("H7").Value = 10
$newSheet.Range("A8").Value = 6
$newSheet.Range("B8").Value = "006648"
$newSheet.Range("C8").Value = "汇安多因子混合A"
$newSheet.Range("D8").Value = "3.69"
$newSheet.Range("E8").Value = "93.49"
$newSheet.Range("F8").Value = "2.87"
$newSheet.Range("G8").Value = "0.1059"
$newSheet.Range("H8").Value = 6
$newSheet.Range("A9").Value = 7
$newSheet.Range("B9").Value = "013867"
$newSheet.Range("C9").Value = "汇安优势企业精选混合A"
$newSheet.Range("D9").Value = "3.54"
$newSheet.Range("E9").Value = "94.26"
$newSheet.Range("F9").Value = "2.99"
$newSheet.Range("G9").Value = "0.1058"
$newSheet.Range("H9").Value = 8
$newSheet.Range("A10").Value = 8
$newSheet.Range("B10").Value = "001550"
$newSheet.Range("C10").Value = "天弘中证医药100指数型发起式 A"
$newSheet.Range("D10").Value = "5.31"
$newSheet.Range("E10").Value = "95.24"
$newSheet.Range("F10").Value = "1.51"
$newSheet.Range("G10").Value = "0.0802"
$newSheet.Range("H10").Value = 2
$newSheet.Range("A11").Value = 9
$newSheet.Range("B11").Value = "160642"
$newSheet.Range("C11").Value = "鹏华增瑞灵活配置混合（LOF）"
$newSheet.Range("D11").Value = "2.05"
$newSheet.Range("E11").Value = "90.96"
$newSheet.Range("F11").Value = "3.13"
$newSheet.Range("G11").Value = "0.0642"
$newSheet.Range("H11").Value = 10
$newSheet.Range("A12").Value = 10
$newSheet.Range("B12").Value = "006649"
$newSheet.Range("C12").Value = "汇安多因子混合C"
$newSheet.Range("D12").Value = "1.95"
$newSheet.Range("E12").Value = "93.49"
$newSheet.Range("F12").Value = "2.87"
$newSheet.Range("G12").Value = "0.0560"
$newSheet.Range("H12").Value = 6
$newSheet.Range("A13").Value = 11
$newSheet.Range("B13").Value = "010421"
$newSheet.Range("C13").Value = "海富通消费优选混合A"
$newSheet.Range("D13").Value = "1.22"
$newSheet.Range("E13").Value = "92.47"
$newSheet.Range("F13").Value = "4.33"
$newSheet.Range("G13").Value = "0.0528"
$newSheet.Range("H13").Value = 9
$newSheet.Range("A14").Value = 12
$newSheet.Range("B14").Value = "002133"
$newSheet.Range("C14").Value = "广发鑫益灵活配置混合"
$newSheet.Range("D14").Value = "1.22"
$newSheet.Range("E14").Value = "93.85"
$newSheet.Range("F14").Value = "3.83"
$newSheet.Range("G14").Value = "0.0467"
$newSheet.Range("H14").Value = 10
$newSheet.Range("A15").Value = 13
$newSheet.Range("B15").Value = "006048"
$newSheet.Range("C15").Value = "长城中证500指数增强A"
$newSheet.Range("D15").Value = "2.12"
$newSheet.Range("E15").Value = "94.78"
$newSheet.Range("F15").Value = "1.80"
$newSheet.Range("G15").Value = "0.0382"
$newSheet.Range("H15").Value = 10
$newSheet.Range("A16").Value = 14
$newSheet.Range("B16").Value = "010558"
$newSheet.Range("C16").Value = "汇安鑫利优选混合A"
$newSheet.Range("D16").Value = "1.30"
$newSheet.Range("E16").Value = "93.99"
$newSheet.Range("F16").Value = "2.87"
$newSheet.Range("G16").Value = "0.0373"
$newSheet.Range("H16").Value = 6
$newSheet.Range("A17").Value = 15
$newSheet.Range("B17").Value = "010221"
$newSheet.Range("C17").Value = "海富通消费核心资产混合C"
$newSheet.Range("D17").Value = "0.57"
$newSheet.Range("E17").Value = "91.44"
$newSheet.Range("F17").Value = "4.25"
$newSheet.Range("G17").Value = "0.0242"
$newSheet.Range("H17").Value = 9
$newSheet.Range("A18").Value = 16
$newSheet.Range("B18").Value = "003889"
$newSheet.Range("C18").Value = "汇安丰泽灵活配置混合A"
$newSheet.Range("D18").Value = "0.64"
$newSheet.Range("E18").Value = "93.72"
$newSheet.Range("F18").Value = "3.19"
$newSheet.Range("G18").Value = "0.0204"
$newSheet.Range("H18").Value = 6
$newSheet.Range("A19").Value = 17
$newSheet.Range("B19").Value = "007413"
$newSheet.Range("C19").Value = "长城中证500指数增强C"
$newSheet.Range("D19").Value = "1.06"
$newSheet.Range("E19").Value = "94.78"
$newSheet.Range("F19").Value = "1.80"
$newSheet.Range("G19").Value = "0.0191"
$newSheet.Range("H19").Value = 10
$newSheet.Range("A20").Value = 18
$newSheet.Range("B20").Value = "010559"
$newSheet.Range("C20").Value = "汇安鑫利优选混合C"
$newSheet.Range("D20").Value = "0.65"
$newSheet.Range("E20").Value = "93.99"
$newSheet.Range("F20").Value = "2.87"
$newSheet.Range("G20").Value = "0.0187"
$newSheet.Range("H20").Value = 6
$newSheet.Range("A21").Value = 19
$newSheet.Range("B21").Value = "010422"
$newSheet.Range("C21").Value = "海富通消费优选混合C"
$newSheet.Range("D21").Value = "0.30"
$newSheet.Range("E21").Value = "92.47"
$newSheet.Range("F21").Value = "4.33"
$newSheet.Range("G21").Value = "0.0130"
$newSheet.Range("H21").Value = 9
$newSheet.Range("A22").Value = 20
$newSheet.Range("B22").Value = "003890"
$newSheet.Range("C22").Value = "汇安丰泽灵活配置混合C"
$newSheet.Range("D22").Value = "0.28"
$newSheet.Range("E22").Value = "93.72"
$newSheet.Range("F22").Value = "3.19"
$newSheet.Range("G22").Value = "0.0089"
$newSheet.Range("H22").Value = 6
$newSheet.Range("A23").Value = 21
$newSheet.Range("B23").Value = "013920"
$newSheet.Range("C23").Value = "兴华创新医疗6个月持有混合A"
$newSheet.Range("D23").Value = "0.18"
$newSheet.Range("E23").Value = "94.83"
$newSheet.Range("F23").Value = "3.99"
$newSheet.Range("G23").Value = "0.0072"
$newSheet.Range("H23").Value = 9
$newSheet.Range("A24").Value = 22
$newSheet.Range("B24").Value = "013868"
$newSheet.Range("C24").Value = "汇安优势企业精选混合C"
$newSheet.Range("D24").Value = "0.23"
$newSheet.Range("E24").Value = "94.26"
$newSheet.Range("F24").Value = "2.99"
$newSheet.Range("G24").Value = "0.0069"
$newSheet.Range("H24").Value = 8
$newSheet.Range("A25").Value = 23
$newSheet.Range("B25").Value = "013921"
$newSheet.Range("C25").Value = "兴华创新医疗6个月持有混合C"
$newSheet.Range("D25").Value = "0.05"
$newSheet.Range("E25").Value = "94.83"
$newSheet.Range("F25").Value = "3.99"
$newSheet.Range("G25").Value = "0.0020"
$newSheet.Range("H25").Value = 9

# --- Step 6: restore "总计" as the active sheet/tab (Worksheets.Add activates the new sheet) ---
$total.Activate()

Write-Host "Edit complete"
